$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 (event "610") is removed entirely; all subsequent rows shift up by one.
$ws.Rows(30).Delete()

# Column A holds numeric-looking event codes stored as TEXT (not numbers).
# Temporarily mark the range as Text so Excel does not auto-convert the
# assigned strings to numbers, then clear the formatting override so the
# cells end up with their original (General) style, same as the rest of the sheet.
$colA = $ws.Range("A2:A38")
$colA.NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "100"
$ws.Range("B2").Value = "Accidente ofidico"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("A3").Value = "113"
$ws.Range("B3").Value = "Desnutrici”n aguda en menores de 5 anos"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.27

# Row 4
$ws.Range("A4").Value = "115"
$ws.Range("B4").Value = "Cancer en menores de 18 anos"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("A5").Value = "155"
$ws.Range("B5").Value = "Cancer de la mama y cuello uterino"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 0.2

# Row 6
$ws.Range("A6").Value = "210"
$ws.Range("B6").Value = "Dengue"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 33
$ws.Range("E6").Value = 0

# Row 7
$ws.Range("A7").Value = "215"
$ws.Range("B7").Value = "Defectos congenitos"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.05

# Row 8
$ws.Range("A8").Value = "217"
$ws.Range("B8").Value = "Chikungunya"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1

# Row 9
$ws.Range("A9").Value = "220"
$ws.Range("B9").Value = "Dengue grave"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1

# Row 10
$ws.Range("A10").Value = "298"
$ws.Range("B10").Value = "Evento adverso grave posterior a la vacunacion"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("A11").Value = "300"
$ws.Range("B11").Value = "Agresiones por animales potencialmente transmisores de rabia"
$ws.Range("C11").Value = 40
$ws.Range("D11").Value = 25
$ws.Range("E11").Value = 0

# Row 12
$ws.Range("A12").Value = "330"
$ws.Range("B12").Value = "Hepatitis a"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1

# Row 13
$ws.Range("A13").Value = "340"
$ws.Range("B13").Value = "Hepatitis b, c y coinfeccion hepatitis b y delta"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.37

# Row 14
$ws.Range("A14").Value = "342"
$ws.Range("B14").Value = "Enfermedades huerfanas - raras"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0.22

# Row 15
$ws.Range("A15").Value = "346"
$ws.Range("B15").Value = "Ira por virus nuevo"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.37

# Row 16
$ws.Range("A16").Value = "348"
$ws.Range("B16").Value = "Infeccion respiratoria aguda grave irag inusitada"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1

# Row 17
$ws.Range("A17").Value = "352"
$ws.Range("B17").Value = "Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0.37

# Row 18
$ws.Range("A18").Value = "355"
$ws.Range("B18").Value = "Enfermedad transmitida por alimentos o agua (eta)"
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.02

# Row 19
$ws.Range("A19").Value = "356"
$ws.Range("B19").Value = "Intento de suicidio"
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0.02

# Row 20
$ws.Range("A20").Value = "357"
$ws.Range("B20").Value = "Iad - infecciones asociadas a dispositivos - individual"
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.14

# Row 21
$ws.Range("A21").Value = "365"
$ws.Range("B21").Value = "Intoxicaciones"
$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0.01

# Row 22
$ws.Range("A22").Value = "420"
$ws.Range("B22").Value = "Leishmaniasis cutanea"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1

# Row 23
$ws.Range("A23").Value = "430"
$ws.Range("D23").Value = 0

# Row 24
$ws.Range("A24").Value = "455"
$ws.Range("B24").Value = "Leptospirosis"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0.37

# Row 25
$ws.Range("A25").Value = "465"
$ws.Range("B25").Value = "Malaria"
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0

# Row 26
$ws.Range("A26").Value = "535"
$ws.Range("B26").Value = "Meningitis bacteriana y enfermedad meningoc”cica"
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1

# Row 27
$ws.Range("A27").Value = "549"
$ws.Range("B27").Value = "Morbilidad materna extrema"
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = 0.1

# Row 28
$ws.Range("A28").Value = "560"
$ws.Range("B28").Value = "Mortalidad perinatal y neonatal tardia"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0

# Row 29
$ws.Range("A29").Value = "580"
$ws.Range("B29").Value = "Mortalidad por dengue"
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1

# Row 30
$ws.Range("A30").Value = "620"
$ws.Range("B30").Value = "Parotiditis"
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0.14

# Row 31
$ws.Range("A31").Value = "720"
$ws.Range("B31").Value = "Sindrome de rubeola congenita"
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 1

# Row 32
$ws.Range("A32").Value = "730"
$ws.Range("B32").Value = "Sarampion"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 1

# Row 33
$ws.Range("A33").Value = "740"
$ws.Range("B33").Value = "Sifilis congenita"
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 1

# Row 34
$ws.Range("A34").Value = "750"
$ws.Range("B34").Value = "Sifilis gestacional"
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 0.04

# Row 35
$ws.Range("A35").Value = "813"
$ws.Range("B35").Value = "Tuberculosis"
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 2
$ws.Range("E35").Value = 0.04

# Row 36
$ws.Range("A36").Value = "831"
$ws.Range("B36").Value = "Varicela individual"
$ws.Range("C36").Value = 10
$ws.Range("D36").Value = 13
$ws.Range("E36").Value = 0.07000000000000001

# Row 37
$ws.Range("A37").Value = "850"
$ws.Range("B37").Value = "Vih/sida/mortalidad por sida"
$ws.Range("C37").Value = 8
$ws.Range("D37").Value = 3
$ws.Range("E37").Value = 0.03

# Row 38
$ws.Range("A38").Value = "895"
$ws.Range("B38").Value = "Zika"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 1

# Restore the original (unstyled/General) formatting for column A now that the
# text values are locked in.
$colA.ClearFormats()